# Auto-generated Excel COM-interop script to apply value updates
# derived from the Omega_Profits.xlsx diff (per-sheet cell value changes).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 251.5
$ws.Range("I9").Value = 10
$ws.Range("K9").Value = 10
$ws.Range("M9").Value = 159
$ws.Range("H19").Value = 2980.647
$ws.Range("I19").Value = 2696.8572
$ws.Range("J19").Value = 3179.3
$ws.Range("K19").Value = 2696.8572
$ws.Range("L19").Value = 3179.3
$ws.Range("M19").Value = -2521.8572
$ws.Range("N19").Value = -3529.3
$ws.Range("H52").Value = 150
$ws.Range("I52").Value = 150
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 450
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -290
$ws.Range("H62").Value = 4953.154
$ws.Range("I62").Value = 4656.143
$ws.Range("K62").Value = 4656.143
$ws.Range("M62").Value = -4032.143
$ws.Range("H65").Value = 4953.154
$ws.Range("I65").Value = 4656.143
$ws.Range("K65").Value = 23280.715
$ws.Range("M65").Value = -20160.715
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
$ws.Range("H100").Value = 6044.5
$ws.Range("I100").Value = 6044.5
$ws.Range("K100").Value = 6044.5
$ws.Range("M100").Value = -5503.5
$ws.Range("H107").Value = 1190.8125
$ws.Range("I107").Value = 931.0769
$ws.Range("K107").Value = 931.0769
$ws.Range("M107").Value = 988.9231
$ws.Range("H113").Value = 1824.4445
$ws.Range("I113").Value = 1824.4445
$ws.Range("K113").Value = 1824.4445
$ws.Range("M113").Value = 1429.5555
$ws.Range("H118").Value = 819.6667
$ws.Range("I118").Value = 768.1429000000001
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 2304.4287
$ws.Range("L118").Value = 3000
$ws.Range("M118").Value = -647.4287000000004
$ws.Range("N118").Value = -6314
$ws.Range("H125").Value = 1623
$ws.Range("J125").Value = 2000
$ws.Range("L125").Value = 18000
$ws.Range("N125").Value = -22920
$ws.Range("H129").Value = 1634.7142
$ws.Range("I129").Value = 1634.7142
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 4904.142599999999
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 95.85740000000078
$ws.Range("H132").Value = 3053.0571
$ws.Range("I132").Value = 2748.7
$ws.Range("K132").Value = 8246.099999999999
$ws.Range("M132").Value = -5716.099999999999
$ws.Range("H137").Value = 2459.8462
$ws.Range("I137").Value = 2555.4443
$ws.Range("J137").Value = 2244.75
$ws.Range("K137").Value = 7666.3329
$ws.Range("L137").Value = 6734.25
$ws.Range("M137").Value = -5116.3329
$ws.Range("N137").Value = -11834.25
$ws.Range("H138").Value = 2635.8965
$ws.Range("I138").Value = 1486.4736
$ws.Range("J138").Value = 3195.8718
$ws.Range("K138").Value = 4459.4208
$ws.Range("L138").Value = 9587.615399999999
$ws.Range("M138").Value = 680.5792000000001
$ws.Range("N138").Value = -19867.6154
$ws.Range("H141").Value = 3090.5862
$ws.Range("I141").Value = 3228.2917
$ws.Range("J141").Value = 2429.6
$ws.Range("K141").Value = 9684.875100000001
$ws.Range("L141").Value = 7288.799999999999
$ws.Range("M141").Value = -4504.875100000001
$ws.Range("N141").Value = -17648.8
$ws.Range("N52").ClearContents()
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1500.8148
$ws.Range("I2").Value = 1052.1
$ws.Range("K2").Value = 1052.1
$ws.Range("M2").Value = -939.0999999999999
$ws.Range("H61").Value = 3571.75
$ws.Range("I61").Value = 3313.9412
$ws.Range("J61").Value = 5032.6665
$ws.Range("K61").Value = 3313.9412
$ws.Range("L61").Value = 5032.6665
$ws.Range("M61").Value = -3101.9412
$ws.Range("N61").Value = -5456.6665
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H74").Value = 2693.8125
$ws.Range("I74").Value = 2406.7334
$ws.Range("K74").Value = 2406.7334
$ws.Range("M74").Value = -1532.7334
$ws.Range("H77").Value = 2693.8125
$ws.Range("I77").Value = 2406.7334
$ws.Range("K77").Value = 12033.667
$ws.Range("M77").Value = -7665.667000000001
$ws.Range("H116").Value = 1500.8148
$ws.Range("I116").Value = 1052.1
$ws.Range("K116").Value = 1052.1
$ws.Range("M116").Value = 1241.9
$ws.Range("H132").Value = 2648.8
$ws.Range("I132").Value = 2526.7222
$ws.Range("K132").Value = 7580.1666
$ws.Range("M132").Value = -5050.1666
$ws.Range("H136").Value = 3571.75
$ws.Range("I136").Value = 3313.9412
$ws.Range("J136").Value = 5032.6665
$ws.Range("K136").Value = 9941.8236
$ws.Range("L136").Value = 15097.9995
$ws.Range("M136").Value = -7391.8236
$ws.Range("N136").Value = -20197.9995
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1500.8148
$ws.Range("I3").Value = 1052.1
$ws.Range("K3").Value = 1052.1
$ws.Range("M3").Value = -938.0999999999999
$ws.Range("H7").Value = 10719983
$ws.Range("I7").Value = 16669583
$ws.Range("J7").Value = 6257783
$ws.Range("K7").Value = 16669583
$ws.Range("L7").Value = 6257783
$ws.Range("M7").Value = -16669470
$ws.Range("N7").Value = -6258009
$ws.Range("H94").Value = 1971.3667
$ws.Range("J94").Value = 2085.5454
$ws.Range("L94").Value = 2085.5454
$ws.Range("N94").Value = -2987.5454
$ws.Range("H105").Value = 3274.05
$ws.Range("I105").Value = 3340.5881
$ws.Range("K105").Value = 3340.5881
$ws.Range("M105").Value = -1593.5881
$ws.Range("H107").Value = 2261.8096
$ws.Range("I107").Value = 2000.375
$ws.Range("J107").Value = 3098.4
$ws.Range("K107").Value = 2000.375
$ws.Range("L107").Value = 3098.4
$ws.Range("M107").Value = -80.375
$ws.Range("N107").Value = -6938.4
$ws.Range("H134").Value = 3072.3462
$ws.Range("I134").Value = 3127.6956
$ws.Range("K134").Value = 9383.086800000001
$ws.Range("M134").Value = -6848.086800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12111.714
$ws.Range("I16").Value = 2050.2222
$ws.Range("K16").Value = 2050.2222
$ws.Range("M16").Value = -1763.2222
$ws.Range("H22").Value = 1082.6666
$ws.Range("I22").Value = 749
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 749
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = -399
$ws.Range("N22").Value = -2450
$ws.Range("H31").Value = 5559.136
$ws.Range("I31").Value = 8868.947
$ws.Range("J31").Value = 3043.68
$ws.Range("K31").Value = 8868.947
$ws.Range("L31").Value = 3043.68
$ws.Range("M31").Value = -8573.947
$ws.Range("N31").Value = -3633.68
$ws.Range("H34").Value = 5559.136
$ws.Range("I34").Value = 8868.947
$ws.Range("J34").Value = 3043.68
$ws.Range("K34").Value = 8868.947
$ws.Range("L34").Value = 3043.68
$ws.Range("M34").Value = -8666.947
$ws.Range("N34").Value = -3447.68
$ws.Range("H58").Value = 4877.263
$ws.Range("I58").Value = 4592.6665
$ws.Range("K58").Value = 4592.6665
$ws.Range("M58").Value = -4389.6665
$ws.Range("H94").Value = 1213.7142
$ws.Range("I94").Value = 1149.25
$ws.Range("K94").Value = 1149.25
$ws.Range("M94").Value = -698.25
$ws.Range("H105").Value = 4548.778
$ws.Range("I105").Value = 4991.2856
$ws.Range("K105").Value = 4991.2856
$ws.Range("M105").Value = -3244.2856
$ws.Range("H107").Value = 992.4286
$ws.Range("I107").Value = 989
$ws.Range("J107").Value = 1013
$ws.Range("K107").Value = 989
$ws.Range("L107").Value = 1013
$ws.Range("M107").Value = 931
$ws.Range("N107").Value = -4853
$ws.Range("H113").Value = 12111.714
$ws.Range("I113").Value = 2050.2222
$ws.Range("K113").Value = 2050.2222
$ws.Range("M113").Value = 119.7777999999998
$ws.Range("H122").Value = 1006666.3
$ws.Range("I122").Value = 1006666.3
$ws.Range("K122").Value = 3019998.9
$ws.Range("M122").Value = -3017548.9
$ws.Range("H132").Value = 8541.793
$ws.Range("I132").Value = 7908.2666
$ws.Range("J132").Value = 9220.571
$ws.Range("K132").Value = 23724.7998
$ws.Range("L132").Value = 27661.713
$ws.Range("M132").Value = -21194.7998
$ws.Range("N132").Value = -32721.713
$ws.Range("H134").Value = 2548.875
$ws.Range("I134").Value = 2399.2
$ws.Range("J134").Value = 2798.3333
$ws.Range("K134").Value = 7197.599999999999
$ws.Range("L134").Value = 8394.999899999999
$ws.Range("M134").Value = -4662.599999999999
$ws.Range("N134").Value = -13464.9999
$ws.Range("H136").Value = 4877.263
$ws.Range("I136").Value = 4592.6665
$ws.Range("K136").Value = 13777.9995
$ws.Range("M136").Value = -11227.9995
$ws.Range("H141").Value = 100999.2
$ws.Range("J141").Value = 100999.2
$ws.Range("L141").Value = 100999.2
$ws.Range("N141").Value = -111359.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1533.091
$ws.Range("J2").Value = 279
$ws.Range("L2").Value = 1674
$ws.Range("N2").Value = -1900
$ws.Range("H23").Value = 308.27274
$ws.Range("J23").Value = 331.1
$ws.Range("L23").Value = 993.3000000000001
$ws.Range("N23").Value = -1463.3
$ws.Range("H70").Value = 16630.8
$ws.Range("I70").Value = 12461.6
$ws.Range("K70").Value = 37384.8
$ws.Range("M70").Value = -37069.8
$ws.Range("H73").Value = 16630.8
$ws.Range("I73").Value = 12461.6
$ws.Range("K73").Value = 37384.8
$ws.Range("M73").Value = -36292.8
$ws.Range("H75").Value = 3850.4
$ws.Range("I75").Value = 3848
$ws.Range("J75").Value = 3851
$ws.Range("K75").Value = 11544
$ws.Range("L75").Value = 11553
$ws.Range("M75").Value = -10546
$ws.Range("N75").Value = -13549
$ws.Range("H78").Value = 3850.4
$ws.Range("I78").Value = 3848
$ws.Range("J78").Value = 3851
$ws.Range("K78").Value = 34632
$ws.Range("L78").Value = 34659
$ws.Range("M78").Value = -29640
$ws.Range("N78").Value = -44643
$ws.Range("H87").Value = 2998.5
$ws.Range("I87").Value = 2998.5
$ws.Range("K87").Value = 8995.5
$ws.Range("M87").Value = -7747.5
$ws.Range("H90").Value = 2998.5
$ws.Range("I90").Value = 2998.5
$ws.Range("K90").Value = 26986.5
$ws.Range("M90").Value = -20746.5
$ws.Range("H94").Value = 15149.625
$ws.Range("I94").Value = 10173.5
$ws.Range("K94").Value = 30520.5
$ws.Range("M94").Value = -29844.5
$ws.Range("H103").Value = 840.1667
$ws.Range("I103").Value = 710.5
$ws.Range("J103").Value = 1099.5
$ws.Range("K103").Value = 2131.5
$ws.Range("L103").Value = 3298.5
$ws.Range("M103").Value = -1252.5
$ws.Range("N103").Value = -5056.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H97").Value = 1115.3334
$ws.Range("I97").Value = 1046.75
$ws.Range("K97").Value = 1046.75
$ws.Range("M97").Value = -550.75
$ws.Range("H107").Value = 640.1111
$ws.Range("I107").Value = 570.125
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 570.125
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1349.875
$ws.Range("N107").Value = -5040
$ws.Range("H113").Value = 6786.8667
$ws.Range("I113").Value = 8574.875
$ws.Range("K113").Value = 8574.875
$ws.Range("M113").Value = -6404.875
$ws.Range("H132").Value = 3931.288
$ws.Range("I132").Value = 3956.4807
$ws.Range("J132").Value = 3744.1428
$ws.Range("K132").Value = 11869.4421
$ws.Range("L132").Value = 11232.4284
$ws.Range("M132").Value = -9339.4421
$ws.Range("N132").Value = -16292.4284
$ws.Range("H133").Value = 82043.57000000001
$ws.Range("J133").Value = 82043.57000000001
$ws.Range("L133").Value = 82043.57000000001
$ws.Range("N133").Value = -92163.57000000001
$ws.Range("H136").Value = 26758.309
$ws.Range("I136").Value = 12323
$ws.Range("J136").Value = 27961.25
$ws.Range("K136").Value = 36969
$ws.Range("L136").Value = 83883.75
$ws.Range("M136").Value = -34419
$ws.Range("N136").Value = -88983.75
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11846.27
$ws.Range("I7").Value = 12032.272
$ws.Range("K7").Value = 12032.272
$ws.Range("M7").Value = -11920.272
$ws.Range("H40").Value = 5373.857
$ws.Range("I40").Value = 1628.3334
$ws.Range("K40").Value = 1628.3334
$ws.Range("M40").Value = -1492.3334
$ws.Range("H46").Value = 5148.2104
$ws.Range("I46").Value = 5555.6875
$ws.Range("K46").Value = 5555.6875
$ws.Range("M46").Value = -5367.6875
$ws.Range("H61").Value = 2547.5833
$ws.Range("I61").Value = 2487.15
$ws.Range("K61").Value = 2487.15
$ws.Range("M61").Value = -2285.15
$ws.Range("H82").Value = 1992.9
$ws.Range("I82").Value = 1376.2
$ws.Range("J82").Value = 2609.6
$ws.Range("K82").Value = 1376.2
$ws.Range("L82").Value = 2609.6
$ws.Range("M82").Value = -1015.2
$ws.Range("N82").Value = -3331.6
$ws.Range("H85").Value = 1992.9
$ws.Range("I85").Value = 1376.2
$ws.Range("J85").Value = 2609.6
$ws.Range("K85").Value = 1376.2
$ws.Range("L85").Value = 2609.6
$ws.Range("M85").Value = -128.2
$ws.Range("N85").Value = -5105.6
$ws.Range("H93").Value = 1939
$ws.Range("I93").Value = 1939
$ws.Range("K93").Value = 1939
$ws.Range("M93").Value = -691
$ws.Range("H113").Value = 2547.5833
$ws.Range("I113").Value = 2487.15
$ws.Range("K113").Value = 2487.15
$ws.Range("M113").Value = -317.1500000000001
$ws.Range("H122").Value = 6736.5884
$ws.Range("I122").Value = 6938.9375
$ws.Range("K122").Value = 20816.8125
$ws.Range("M122").Value = -18366.8125
$ws.Range("H126").Value = 11846.27
$ws.Range("I126").Value = 12032.272
$ws.Range("K126").Value = 36096.81600000001
$ws.Range("M126").Value = -33626.81600000001
$ws.Range("H132").Value = 2358.92
$ws.Range("I132").Value = 1894.1333
$ws.Range("K132").Value = 5682.3999
$ws.Range("M132").Value = -3152.3999
$ws.Range("H136").Value = 1967.3334
$ws.Range("I136").Value = 1666
$ws.Range("K136").Value = 4998
$ws.Range("M136").Value = -2448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10123.75
$ws.Range("I62").Value = 16375
$ws.Range("J62").Value = 3872.5
$ws.Range("K62").Value = 16375
$ws.Range("L62").Value = 3872.5
$ws.Range("M62").Value = -15751
$ws.Range("N62").Value = -5120.5
$ws.Range("H65").Value = 10123.75
$ws.Range("I65").Value = 16375
$ws.Range("J65").Value = 3872.5
$ws.Range("K65").Value = 81875
$ws.Range("L65").Value = 19362.5
$ws.Range("M65").Value = -78755
$ws.Range("N65").Value = -25602.5
$ws.Range("H100").Value = 2328.3333
$ws.Range("I100").Value = 2335.7144
$ws.Range("J100").Value = 2318
$ws.Range("K100").Value = 4671.4288
$ws.Range("L100").Value = 4636
$ws.Range("M100").Value = -4130.4288
$ws.Range("N100").Value = -5718
$ws.Range("H122").Value = 3704
$ws.Range("I122").Value = 3538.6667
$ws.Range("K122").Value = 10616.0001
$ws.Range("M122").Value = -8166.000100000001
$ws.Range("H126").Value = 2348.2
$ws.Range("I126").Value = 2310.5
$ws.Range("K126").Value = 6931.5
$ws.Range("M126").Value = -4461.5
$ws.Range("H132").Value = 2576.587
$ws.Range("I132").Value = 2434.756
$ws.Range("K132").Value = 7304.268
$ws.Range("M132").Value = -4774.268
$ws.Range("H136").Value = 4626.769
$ws.Range("I136").Value = 3941.7273
$ws.Range("K136").Value = 11825.1819
$ws.Range("M136").Value = -9275.1819
